$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.333.47"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "3.667.72"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "644.03"
$ws.Range("E5").Value = "  -5.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.35"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.496"
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("E9").Value = "  -0.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.08"
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.442"
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000232"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "4.284.92"
$ws.Range("E13").Value = "  -0.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.53"
$ws.Range("E14").Value = "  +0.11%  "
$ws.Range("D15").Value = "3.665.48"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").Value = "69.317.42"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "16.02"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.46"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "465.13"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.87"
$ws.Range("E21").Value = "  -0.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.644"
$ws.Range("E22").Value = "  -1.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.43"
$ws.Range("E23").Value = "  -0.69%  "
$ws.Range("D24").Value = "3.813.63"
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000125"
$ws.Range("E26").Value = "  +1.32%  "
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.99"
$ws.Range("E28").Value = "  -1.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.62"
$ws.Range("E29").Value = "  -2.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.71"
$ws.Range("E30").Value = "  -1.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.00"
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.69"
$ws.Range("E33").Value = "  -0.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.45"
$ws.Range("E34").Value = "  -2.78%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.163"
$ws.Range("E35").Value = "  +4.40%  "
$ws.Range("D36").Value = "3.656.59"
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.42"
$ws.Range("E37").Value = "  +1.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.90"
$ws.Range("E39").Value = "  -5.95%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "177.45"
$ws.Range("E41").Value = "  +4.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0899"
$ws.Range("E42").Value = "  -0.79%  "
$ws.Range("E43").Value = "  -3.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.923"
$ws.Range("E44").Value = "  -2.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.65"
$ws.Range("E45").Value = "  -2.09%  "
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.05"
$ws.Range("E47").Value = "  -4.26%  "
$ws.Range("B48").Value = "FLOKI"
$ws.Range("C48").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000268"
$ws.Range("E48").Value = "  -2.94%  "
$ws.Range("B49").Value = "SuiNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.07"
$ws.Range("E49").Value = "  -4.58%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.25"
$ws.Range("E50").Value = "  -3.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.83"
$ws.Range("E51").Value = "  +0.47%  "
